$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.008.15'
$ws.Range("E2").Value = '  +0.21%  '

$ws.Range("D3").Value = '3.519.42'
$ws.Range("E3").Value = '  +0.07%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.66%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.64'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.86%  '

$ws.Range("D7").Value = '3.518.07'
$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("E9").Value = '  +0.81%  '

$ws.Range("E10").Value = '  +1.28%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.15'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.55%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.384'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.19%  '

$ws.Range("D13").Value = '4.113.52'
$ws.Range("E13").Value = '  +0.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.23'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.70%  '

$ws.Range("E15").Value = '  +0.69%  '

$ws.Range("E16").Value = '  +0.04%  '

$ws.Range("D17").Value = '3.521.76'
$ws.Range("E17").Value = '  -0.18%  '

$ws.Range("D18").Value = '64.998.28'
$ws.Range("E18").Value = '  +0.17%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.84'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.79%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.42'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.04%  '

$ws.Range("E21").Value = '  -1.58%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '389.67'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.33%  '

$ws.Range("E23").Value = '  +1.39%  '

$ws.Range("D24").Value = '3.661.19'
$ws.Range("E24").Value = '  -0.09%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.25'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.83%  '

$ws.Range("E27").Value = '  +2.15%  '

$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.81'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.71%  '

$ws.Range("B29").Value = 'Fetch.AI'
$ws.Range("C29").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.61'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +19.70%  '

$ws.Range("E30").Value = '  +0.50%  '

$ws.Range("E31").Value = '  +2.03%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.40'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.96%  '

$ws.Range("D33").Value = '3.525.35'
$ws.Range("E33").Value = '  -0.19%  '

$ws.Range("E34").Value = '  +1.80%  '

$ws.Range("E35").Value = '  +0.04%  '

$ws.Range("E36").Value = '  +2.34%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.29'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.50%  '

$ws.Range("E38").Value = '  +1.87%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '169.84'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.53%  '

$ws.Range("E40").Value = '  +0.40%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0825'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.65%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.821'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.53%  '

$ws.Range("B43").Value = 'ONDO'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.24'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.13%  '

$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.59'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.68%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.08%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.26'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.35%  '

$ws.Range("E47").Value = '  +1.11%  '

$ws.Range("E48").Value = '  +0.43%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.94'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.77%  '

$ws.Range("D50").Value = '2.359.51'
$ws.Range("E50").Value = '  -0.97%  '

$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0268'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.50%  '
